$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.162.62'
$ws.Range('E2').Value = '  +0.30%  '
$ws.Range('D3').Value = '1.905.46'
$ws.Range('E3').Value = '  +0.77%  '
$ws.Range('D4').Value = "'1.001"
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = "'306.42"
$ws.Range('E5').Value = '  -0.18%  '
$ws.Range('D7').Value = "'0.5233"
$ws.Range('D8').Value = "'0.3771"
$ws.Range('E8').Value = '  +0.49%  '
$ws.Range('E9').Value = '  +0.53%  '
$ws.Range('D10').Value = "'21.20"
$ws.Range('E10').Value = '  -0.01%  '
$ws.Range('D11').Value = "'0.9040"
$ws.Range('E11').Value = '  -0.20%  '
$ws.Range('D12').Value = "'0.08519"
$ws.Range('E12').Value = '  +11.45%  '
$ws.Range('B13').Value = 'Litecoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D13').Value = "'96.75"
$ws.Range('E13').Value = '  +1.90%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.901.25'
$ws.Range('E14').Value = '  +0.58%  '
$ws.Range('E15').Value = '  +0.42%  '
$ws.Range('E16').Value = '  +0.14%  '
$ws.Range('D17').Value = "'0.000008643"
$ws.Range('E17').Value = '  +2.01%  '
$ws.Range('D18').Value = "'14.55"
$ws.Range('E18').Value = '  +0.60%  '
$ws.Range('E19').Value = '  +0.13%  '
$ws.Range('D20').Value = '27.197.00'
$ws.Range('E20').Value = '  +0.32%  '
$ws.Range('E21').Value = '  +0.08%  '
$ws.Range('D22').Value = '2.152.95'
$ws.Range('E22').Value = '  +1.34%  '
$ws.Range('E23').Value = '  +0.45%  '
$ws.Range('D24').Value = "'6.437"
$ws.Range('E24').Value = '  +0.52%  '
$ws.Range('D25').Value = "'2.315"
$ws.Range('E25').Value = '  +1.55%  '
$ws.Range('D26').Value = "'147.10"
$ws.Range('E26').Value = '  +0.78%  '
$ws.Range('D27').Value = "'18.25"
$ws.Range('E27').Value = '  +0.97%  '
$ws.Range('D28').Value = "'1.747"
$ws.Range('E28').Value = '  -1.13%  '
$ws.Range('D29').Value = "'115.11"
$ws.Range('E29').Value = '  +0.55%  '
$ws.Range('D30').Value = "'4.929"
$ws.Range('E30').Value = '  -0.45%  '
$ws.Range('D31').Value = "'4.816"
$ws.Range('E31').Value = '  -0.35%  '
$ws.Range('E32').Value = '  +1.40%  '
$ws.Range('D33').Value = "'0.8054"
$ws.Range('E33').Value = '  +2.80%  '
$ws.Range('D34').Value = "'0.05059"
$ws.Range('E34').Value = '  -0.55%  '
$ws.Range('D35').Value = "'1.244"
$ws.Range('E35').Value = '  +0.49%  '
$ws.Range('D36').Value = "'3.450"
$ws.Range('E36').Value = '  +4.95%  '
$ws.Range('D37').Value = "'2.950"
$ws.Range('E37').Value = '  -1.03%  '
$ws.Range('D38').Value = "'2.613"
$ws.Range('E38').Value = '  -0.36%  '
$ws.Range('D39').Value = "'0.5719"
$ws.Range('E39').Value = '  +2.23%  '
$ws.Range('D40').Value = "'0.02001"
$ws.Range('E40').Value = '  -0.01%  '
$ws.Range('D42').Value = "'9.142"
$ws.Range('E42').Value = '  +0.25%  '
$ws.Range('D43').Value = "'6.640"
$ws.Range('E43').Value = '  +0.03%  '
$ws.Range('D44').Value = "'116.26"
$ws.Range('D45').Value = "'0.1517"
$ws.Range('E45').Value = '  +0.57%  '
$ws.Range('D46').Value = "'0.4863"
$ws.Range('E46').Value = '  +1.22%  '
$ws.Range('D47').Value = "'10.17"
$ws.Range('E47').Value = '  -0.40%  '
$ws.Range('D48').Value = "'0.9999"
$ws.Range('E48').Value = '  +0.07%  '
$ws.Range('E49').Value = '  +1.17%  '
$ws.Range('E50').Value = '  +0.18%  '
$ws.Range('D51').Value = "'64.19"
$ws.Range('E51').Value = '  +0.28%  '
